$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp column (O) for every data row (2-390) to reflect the new crawl time.
for ($r = 2; $r -le 390; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-03-11 12:54:44"
}

# Row 188: the product aria label text changed to note it is out of stock online.
$ws.Range("M188").Value = "Kale Dürüm 30cm 18 Stück - Online kein Bestand 40% Aktion 4.50 Schweizer Franken statt 7.50 Schweizer Franken"
